$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 3342.1428
$ws.Cells.Item(74, 9).Value = 3520
$ws.Cells.Item(74, 10).Value = 2897.5
$ws.Cells.Item(74, 11).Value = 3520
$ws.Cells.Item(74, 12).Value = 2897.5
$ws.Cells.Item(74, 13).Value = -2584
$ws.Cells.Item(74, 14).Value = -4769.5
$ws.Cells.Item(77, 8).Value = 3342.1428
$ws.Cells.Item(77, 9).Value = 3520
$ws.Cells.Item(77, 10).Value = 2897.5
$ws.Cells.Item(77, 11).Value = 17600
$ws.Cells.Item(77, 12).Value = 14487.5
$ws.Cells.Item(77, 13).Value = -12920
$ws.Cells.Item(77, 14).Value = -23847.5
$ws.Cells.Item(115, 8).Value = 2076.8
$ws.Cells.Item(115, 9).Value = 2128.1667
$ws.Cells.Item(115, 10).Value = 1999.75
$ws.Cells.Item(115, 11).Value = 6384.500100000001
$ws.Cells.Item(115, 12).Value = 5999.25
$ws.Cells.Item(115, 13).Value = -4817.500100000001
$ws.Cells.Item(115, 14).Value = -9133.25
$ws.Cells.Item(117, 8).Value = 48542
$ws.Cells.Item(117, 10).Value = 48542
$ws.Cells.Item(117, 12).Value = 48542
$ws.Cells.Item(117, 14).Value = -57720
$ws.Cells.Item(129, 8).Value = 2977
$ws.Cells.Item(129, 10).Value = 2698
$ws.Cells.Item(129, 12).Value = 8094
$ws.Cells.Item(129, 14).Value = -18094

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(80, 8).Value = 51398.145
$ws.Cells.Item(80, 10).Value = 51398.145
$ws.Cells.Item(80, 12).Value = 51398.145
$ws.Cells.Item(80, 14).Value = -53394.145
$ws.Cells.Item(83, 8).Value = 51398.145
$ws.Cells.Item(83, 10).Value = 51398.145
$ws.Cells.Item(83, 12).Value = 154194.435
$ws.Cells.Item(83, 14).Value = -164178.435
$ws.Cells.Item(104, 8).Value = 41088.25
$ws.Cells.Item(104, 10).Value = 41088.25
$ws.Cells.Item(104, 12).Value = 41088.25
$ws.Cells.Item(104, 14).Value = -48076.25
$ws.Cells.Item(105, 8).Value = 47942.4
$ws.Cells.Item(105, 10).Value = 47942.4
$ws.Cells.Item(105, 12).Value = 47942.4
$ws.Cells.Item(105, 14).Value = -54930.4
$ws.Cells.Item(107, 8).Value = 36254
$ws.Cells.Item(107, 10).Value = 36254
$ws.Cells.Item(107, 12).Value = 36254
$ws.Cells.Item(107, 14).Value = -43934
$ws.Cells.Item(109, 8).Value = 43149.4
$ws.Cells.Item(109, 10).Value = 43149.4
$ws.Cells.Item(109, 12).Value = 43149.4
$ws.Cells.Item(109, 14).Value = -45923.4
$ws.Cells.Item(113, 8).Value = 46342
$ws.Cells.Item(113, 10).Value = 46342
$ws.Cells.Item(113, 12).Value = 46342
$ws.Cells.Item(113, 14).Value = -55020
$ws.Cells.Item(117, 8).Value = 47910.168
$ws.Cells.Item(117, 10).Value = 47910.168
$ws.Cells.Item(117, 12).Value = 47910.168
$ws.Cells.Item(117, 14).Value = -57088.168
$ws.Cells.Item(118, 8).Value = 49358.332
$ws.Cells.Item(118, 10).Value = 49358.332
$ws.Cells.Item(118, 12).Value = 49358.332
$ws.Cells.Item(118, 14).Value = -52672.332
$ws.Cells.Item(119, 8).Value = 52579.5
$ws.Cells.Item(119, 10).Value = 52579.5
$ws.Cells.Item(119, 12).Value = 52579.5
$ws.Cells.Item(119, 14).Value = -62255.5
$ws.Cells.Item(130, 8).Value = 38187.5
$ws.Cells.Item(130, 10).Value = 38187.5
$ws.Cells.Item(130, 12).Value = 38187.5
$ws.Cells.Item(130, 14).Value = -48227.5
$ws.Cells.Item(132, 8).Value = 16668279
$ws.Cells.Item(132, 9).Value = 31251038
$ws.Cells.Item(132, 10).Value = 2268.5715
$ws.Cells.Item(132, 11).Value = 93753114
$ws.Cells.Item(132, 12).Value = 6805.7145
$ws.Cells.Item(132, 13).Value = -93750584
$ws.Cells.Item(132, 14).Value = -11865.7145
$ws.Cells.Item(134, 8).Value = 44866.273
$ws.Cells.Item(134, 10).Value = 44866.273
$ws.Cells.Item(134, 12).Value = 44866.273
$ws.Cells.Item(134, 14).Value = -55006.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(109, 8).Value = 28179.3
$ws.Cells.Item(109, 10).Value = 28179.3
$ws.Cells.Item(109, 12).Value = 28179.3
$ws.Cells.Item(109, 14).Value = -30259.3
$ws.Cells.Item(111, 8).Value = 47258.332
$ws.Cells.Item(111, 10).Value = 47258.332
$ws.Cells.Item(111, 12).Value = 47258.332
$ws.Cells.Item(111, 14).Value = -55438.332
$ws.Cells.Item(116, 8).Value = 47772.5
$ws.Cells.Item(116, 10).Value = 47772.5
$ws.Cells.Item(116, 12).Value = 47772.5
$ws.Cells.Item(116, 14).Value = -56950.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 2597.7114
$ws.Cells.Item(113, 9).Value = 4361.926
$ws.Cells.Item(113, 10).Value = 692.36
$ws.Cells.Item(113, 11).Value = 13085.778
$ws.Cells.Item(113, 12).Value = 2077.08
$ws.Cells.Item(113, 13).Value = -10915.778
$ws.Cells.Item(113, 14).Value = -6417.08
$ws.Cells.Item(122, 8).Value = 7533.125
$ws.Cells.Item(122, 9).Value = 584.4
$ws.Cells.Item(122, 10).Value = 10691.637
$ws.Cells.Item(122, 11).Value = 5259.599999999999
$ws.Cells.Item(122, 12).Value = 96224.73300000001
$ws.Cells.Item(122, 13).Value = -2809.599999999999
$ws.Cells.Item(122, 14).Value = -101124.733
$ws.Cells.Item(131, 8).Value = 4565.5483
$ws.Cells.Item(131, 9).Value = 7195.467
$ws.Cells.Item(131, 10).Value = 2100
$ws.Cells.Item(131, 11).Value = 21586.401
$ws.Cells.Item(131, 12).Value = 6300
$ws.Cells.Item(131, 13).Value = -16546.401
$ws.Cells.Item(131, 14).Value = -16380
$ws.Cells.Item(134, 8).Value = 50005480
$ws.Cells.Item(134, 9).Value = 90913310
$ws.Cells.Item(134, 10).Value = 7019
$ws.Cells.Item(134, 11).Value = 272739930
$ws.Cells.Item(134, 12).Value = 21057
$ws.Cells.Item(134, 13).Value = -272734860
$ws.Cells.Item(134, 14).Value = -31197

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 338453.66
$ws.Cells.Item(80, 9).Value = 505360.5
$ws.Cells.Item(80, 11).Value = 505360.5
$ws.Cells.Item(80, 13).Value = -504362.5
$ws.Cells.Item(83, 8).Value = 338453.66
$ws.Cells.Item(83, 9).Value = 505360.5
$ws.Cells.Item(83, 11).Value = 2526802.5
$ws.Cells.Item(83, 13).Value = -2521810.5
$ws.Cells.Item(97, 8).Value = 5619.636
$ws.Cells.Item(97, 9).Value = 4645
$ws.Cells.Item(97, 11).Value = 4645
$ws.Cells.Item(97, 13).Value = -4149
$ws.Cells.Item(102, 8).Value = 1945.5834
$ws.Cells.Item(102, 9).Value = 2227.4443
$ws.Cells.Item(102, 10).Value = 1100
$ws.Cells.Item(102, 11).Value = 2227.4443
$ws.Cells.Item(102, 12).Value = 1100
$ws.Cells.Item(102, 13).Value = -605.4443000000001
$ws.Cells.Item(102, 14).Value = -4344
$ws.Cells.Item(104, 8).Value = 44916
$ws.Cells.Item(104, 10).Value = 44916
$ws.Cells.Item(104, 12).Value = 44916
$ws.Cells.Item(104, 14).Value = -51904
$ws.Cells.Item(105, 8).Value = 42937.668
$ws.Cells.Item(105, 10).Value = 42937.668
$ws.Cells.Item(105, 12).Value = 42937.668
$ws.Cells.Item(105, 14).Value = -49925.668
$ws.Cells.Item(116, 8).Value = 38936.715
$ws.Cells.Item(116, 10).Value = 38936.715
$ws.Cells.Item(116, 12).Value = 38936.715
$ws.Cells.Item(116, 14).Value = -48114.715
$ws.Cells.Item(118, 8).Value = 38192.668
$ws.Cells.Item(118, 10).Value = 38192.668
$ws.Cells.Item(118, 12).Value = 38192.668
$ws.Cells.Item(118, 14).Value = -41506.668
$ws.Cells.Item(130, 8).Value = 44334.8
$ws.Cells.Item(130, 10).Value = 44334.8
$ws.Cells.Item(130, 12).Value = 44334.8
$ws.Cells.Item(130, 14).Value = -54374.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1600.5
$ws.Cells.Item(22, 9).Value = 500
$ws.Cells.Item(22, 11).Value = 500
$ws.Cells.Item(22, 13).Value = -205
$ws.Cells.Item(27, 8).Value = 1600.5
$ws.Cells.Item(27, 9).Value = 500
$ws.Cells.Item(27, 11).Value = 500
$ws.Cells.Item(27, 13).Value = -393
$ws.Cells.Item(40, 8).Value = 2569.8333
$ws.Cells.Item(40, 9).Value = 2075.9412
$ws.Cells.Item(40, 11).Value = 2075.9412
$ws.Cells.Item(40, 13).Value = -1939.9412
$ws.Cells.Item(82, 8).Value = 9260268
$ws.Cells.Item(82, 9).Value = 1180.5
$ws.Cells.Item(82, 10).Value = 27778442
$ws.Cells.Item(82, 11).Value = 1180.5
$ws.Cells.Item(82, 12).Value = 27778442
$ws.Cells.Item(82, 13).Value = -819.5
$ws.Cells.Item(82, 14).Value = -27779164
$ws.Cells.Item(85, 8).Value = 9260268
$ws.Cells.Item(85, 9).Value = 1180.5
$ws.Cells.Item(85, 10).Value = 27778442
$ws.Cells.Item(85, 11).Value = 1180.5
$ws.Cells.Item(85, 12).Value = 27778442
$ws.Cells.Item(85, 13).Value = 67.5
$ws.Cells.Item(85, 14).Value = -27780938
$ws.Cells.Item(110, 8).Value = 45527
$ws.Cells.Item(110, 10).Value = 45527
$ws.Cells.Item(110, 12).Value = 45527
$ws.Cells.Item(110, 14).Value = -53707
$ws.Cells.Item(119, 8).Value = 47408
$ws.Cells.Item(119, 10).Value = 47408
$ws.Cells.Item(119, 12).Value = 47408
$ws.Cells.Item(119, 14).Value = -57084
$ws.Cells.Item(121, 8).Value = 41206
$ws.Cells.Item(121, 10).Value = 41206
$ws.Cells.Item(121, 12).Value = 41206
$ws.Cells.Item(121, 14).Value = -44700

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 43175.2
$ws.Cells.Item(103, 10).Value = 43175.2
$ws.Cells.Item(103, 12).Value = 43175.2
$ws.Cells.Item(103, 14).Value = -45519.2
$ws.Cells.Item(105, 8).Value = 50052
$ws.Cells.Item(105, 10).Value = 50052
$ws.Cells.Item(105, 12).Value = 50052
$ws.Cells.Item(105, 14).Value = -57040
